$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 25

function Set-TextCell($r, $c, $text) {
    # Force a genuine "text" cell even when the content would otherwise
    # be auto-detected (dates) or be empty, then strip the style that the
    # quote-prefix / number-format trick leaves behind so the cell ends
    # up with the default style, matching a freshly authored row.
    $cell = $ws.Cells.Item($r, $c)
    $cell.NumberFormat = "@"
    if ($text -eq "") {
        $cell.Value = "'"
    } else {
        $cell.Value = $text
    }
    $cell.Style = "Normal"
}

$ws.Cells.Item($row, 1).Value = 111951526
$ws.Cells.Item($row, 2).Value = 103249
Set-TextCell $row 3 "Ovaliderad"
Set-TextCell $row 4 "EN"
$ws.Cells.Item($row, 5).Value = 340
Set-TextCell $row 6 "Ryl"
Set-TextCell $row 7 "Chimaphila umbellata"
Set-TextCell $row 8 "(L.) W. P. C. Barton"
Set-TextCell $row 9 ""
Set-TextCell $row 10 ""
Set-TextCell $row 11 "fullt utvecklade blad"
Set-TextCell $row 12 ""
Set-TextCell $row 14 ""
Set-TextCell $row 16 "Bränsle SO, Öl"
$ws.Cells.Item($row, 17).Value = 619465
$ws.Cells.Item($row, 18).Value = 6344010
$ws.Cells.Item($row, 19).Value = 25
Set-TextCell $row 20 "Kalmar"
Set-TextCell $row 21 "Borgholm"
Set-TextCell $row 22 "Öland"
Set-TextCell $row 23 "Böda"
Set-TextCell $row 24 "Hö-Bor-8992"
Set-TextCell $row 25 "2023-09-07"
Set-TextCell $row 27 "2023-09-07"
Set-TextCell $row 29 "Lokalen såg ganska ok ut men nära vägen fanns stora bestånd av örnbräken som konkurrerar."
$ws.Cells.Item($row, 30).Value = $true
$ws.Cells.Item($row, 31).Value = $false
Set-TextCell $row 32 ""
$ws.Cells.Item($row, 33).Value = $false
Set-TextCell $row 46 ""
Set-TextCell $row 49 "Jan Yngve Andersson"
Set-TextCell $row 50 "Jan Yngve Andersson"
Set-TextCell $row 51 "Floraväkteri Sverige"
